$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "27.053.26"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").Value = "1.798.32"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.71"
$ws.Range("E5").Value = "  -2.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4208"
$ws.Range("E7").Value = "  -2.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3598"
$ws.Range("E8").Value = "  -2.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07143"
$ws.Range("E9").Value = "  -2.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8445"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.23"
$ws.Range("E11").Value = "  -3.57%  "
$ws.Range("D12").Value = "1.880.00"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.297"
$ws.Range("E13").Value = "  -3.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.372"
$ws.Range("E14").Value = "  -3.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06762"
$ws.Range("E15").Value = "  -2.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.18"
$ws.Range("E17").Value = "  -1.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008699"
$ws.Range("E18").Value = "  -3.97%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("E20").Value = "  -3.69%  "
$ws.Range("D21").Value = "27.177.22"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.066"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.01"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "2.053.00"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.926"
$ws.Range("E25").Value = "  -2.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.09"
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.13"
$ws.Range("E27").Value = "  -4.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.033"
$ws.Range("E28").Value = "  -5.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.44"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.653"
$ws.Range("E30").Value = "  -12.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08989"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7262"
$ws.Range("E32").Value = "  -7.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.860"
$ws.Range("E33").Value = "  -4.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.335"
$ws.Range("E34").Value = "  -5.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.098"
$ws.Range("E35").Value = "  -5.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.079"
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05154"
$ws.Range("E38").Value = "  -5.21%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01902"
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1632"
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4978"
$ws.Range("E41").Value = "  -3.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.607"
$ws.Range("E42").Value = "  -7.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.076"
$ws.Range("E43").Value = "  -6.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.962"
$ws.Range("E44").Value = "  -11.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "105.39"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.17"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06298"
$ws.Range("E48").Value = "  -3.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4539"
$ws.Range("E49").Value = "  -5.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.604"
$ws.Range("E50").Value = "  -3.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.719"
$ws.Range("E51").Value = "  -6.91%  "
